$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Actor:" value for the Stock Disponivel use case
$ws.Range("C3").Value = "Funcionário"

# Match the selection recorded in the saved workbook
$ws.Range("C3:D3").Select()
